# Reconciliation report update: refresh the reconciled/bias figures on the
# "Statistics" sheet (columns D/E and M/N, rows 2-15) and the derived
# m/s + Km/h speed stats on the "Speeds" sheet (columns B/C and H/I,
# rows 2-15) to reflect the new scaling pass.

$wb = $excel.ActiveWorkbook

# --- Statistics sheet: Distances (D/E) and Parciais (M/N) ---------------
$statistics = $wb.Worksheets.Item("Statistics")

$statistics.Range("D2").Value = 24.02015624999999
$statistics.Range("E2").Value = 1.1300329124999742
$statistics.Range("M2").Value = 250.15039088411143
$statistics.Range("N2").Value = 0.29934647992362784

$statistics.Range("D3").Value = 8.031721500000014
$statistics.Range("E3").Value = -3.925212439999978
$statistics.Range("M3").Value = 83.57307479169685
$statistics.Range("N3").Value = 0.08671982108812415

$statistics.Range("D4").Value = 10.241316000000003
$statistics.Range("E4").Value = 0.4925409000000105
$statistics.Range("M4").Value = 109.13106377294746
$statistics.Range("N4").Value = -1.3612458571391954

$statistics.Range("D5").Value = 7.503612499999992
$statistics.Range("E5").Value = -5.020897925000009
$statistics.Range("M5").Value = 73.62914792002994
$statistics.Range("N5").Value = 0.15659575999293907

$statistics.Range("D6").Value = 41.679207
$statistics.Range("E6").Value = -3.778393644999987
$statistics.Range("M6").Value = 33.935730261306304
$statistics.Range("N6").Value = 4.346407223666102

$statistics.Range("D7").Value = 5.256855249999999
$statistics.Range("E7").Value = 0.3162955049999985
$statistics.Range("M7").Value = 45.30763442388707
$statistics.Range("N7").Value = 1.6817722487649434

$statistics.Range("D8").Value = 15.226627750000002
$statistics.Range("E8").Value = 1.5941053225000044
$statistics.Range("M8").Value = 158.97508123642908
$statistics.Range("N8").Value = -0.4165636098426546

$statistics.Range("D9").Value = 6.249899499999999
$statistics.Range("E9").Value = 0.3719126924999996
$statistics.Range("M9").Value = 87.64042093754689
$statistics.Range("N9").Value = -0.5515844070551879

$statistics.Range("D10").Value = 45.03974074999999
$statistics.Range("E10").Value = 3.266254130000007
$statistics.Range("M10").Value = 75.00178929657703
$statistics.Range("N10").Value = 0.27854465147656526

$statistics.Range("D11").Value = 18.959232249999996
$statistics.Range("E11").Value = 1.9918035649999943
$statistics.Range("M11").Value = 163.78093162494042
$statistics.Range("N11").Value = 0.008858195877735398

$statistics.Range("D12").Value = 4.0894685
$statistics.Range("E12").Value = -0.4761868424999971
$statistics.Range("M12").Value = 26.648822360846676
$statistics.Range("N12").Value = 0.41756160757423544

$statistics.Range("D13").Value = 13.502141249999998
$statistics.Range("E13").Value = 1.318360969999997
$statistics.Range("M13").Value = 153.46252723917135
$statistics.Range("N13").Value = -0.7757422087534849

$statistics.Range("D14").Value = 22.4719305
$statistics.Range("E14").Value = 2.6068666675000074
$statistics.Range("M14").Value = 161.93272048892686
$statistics.Range("N14").Value = 0.46394988318041896

$statistics.Range("D15").Value = 222.27190900000002
$statistics.Range("E15").Value = -0.1125181875000294
$statistics.Range("M15").Value = 1423.169335238417
$statistics.Range("N15").Value = 4.63461978875398

# --- Speeds sheet: m/s (B/C) and Km/h (H/I) ------------------------------
$speeds = $wb.Worksheets.Item("Speeds")

$speeds.Range("B2").Value = 10.414186663923617
$speeds.Range("C2").Value = 0.07803939839340147
$speeds.Range("H2").Value = 37.49107199012502
$speeds.Range("I2").Value = 0.2809418342162453

$speeds.Range("B3").Value = 10.405375085739303
$speeds.Range("C3").Value = 0.6556751020613024
$speeds.Range("H3").Value = 37.45935030866149
$speeds.Range("I3").Value = 2.3604303674206886

$speeds.Range("B4").Value = 10.655960989090408
$speeds.Range("C4").Value = 0.1272605155492091
$speeds.Range("H4").Value = 38.36145956072547
$speeds.Range("I4").Value = 0.45813785597715284

$speeds.Range("B5").Value = 9.812493371696634
$speeds.Range("C5").Value = 1.1893968345811738
$speeds.Range("H5").Value = 35.324976138107886
$speeds.Range("I5").Value = 4.281828604492226

$speeds.Range("B6").Value = 0.8142124743713647
$speeds.Range("C6").Value = 0.030212199628910322
$speeds.Range("H6").Value = 2.931164907736913
$speeds.Range("I6").Value = 0.10876391866407716

$speeds.Range("B7").Value = 8.61877154100583
$speeds.Range("C7").Value = 0.0605639651641692
$speeds.Range("H7").Value = 31.027577547620986
$speeds.Range("I7").Value = 0.21803027459100913

$speeds.Range("B8").Value = 10.440596818059669
$speeds.Range("C8").Value = 0.09759698963529764
$speeds.Range("H8").Value = 37.58614854501481
$speeds.Range("I8").Value = 0.3513491626870715

$speeds.Range("B9").Value = 14.022692834908291
$speeds.Range("C9").Value = 0.18472466781855088
$speeds.Range("H9").Value = 50.48169420566985
$speeds.Range("I9").Value = 0.6650088041467832

$speeds.Range("B10").Value = 1.6652358128099802
$speeds.Range("C10").Value = 0.01655996400472885
$speeds.Range("H10").Value = 5.994848926115929
$speeds.Range("I10").Value = 0.05961587041702386

$speeds.Range("B11").Value = 8.638584593790208
$speeds.Range("C11").Value = 0.11742826275302161
$speeds.Range("H11").Value = 31.09890453764475
$speeds.Range("I11").Value = 0.4227417459108778

$speeds.Range("B12").Value = 6.51645130922189
$speeds.Range("C12").Value = 0.33973860967753566
$speeds.Range("H12").Value = 23.459224713198804
$speeds.Range("I12").Value = 1.2230589948391284

$speeds.Range("B13").Value = 11.365791869431924
$speeds.Range("C13").Value = 0.17439041783941456
$speeds.Range("H13").Value = 40.91685072995492
$speeds.Range("I13").Value = 0.6278055042218924

$speeds.Range("B14").Value = 7.205999524114177
$speeds.Range("C14").Value = 0.07549829286370431
$speeds.Range("H14").Value = 25.94159828681104
$speeds.Range("I14").Value = 0.2717938543093355

$speeds.Range("B15").Value = 6.402830396523102
$speeds.Range("C15").Value = 0.050528616619386534
$speeds.Range("H15").Value = 23.05018942748317
$speeds.Range("I15").Value = 0.1819030198297915
